$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns with latest scraped values.
# Numeric-looking price strings must be forced back to literal text (matching the source
# data's inline-string formatting) by briefly switching the cell to Text format and then
# restoring the default "Normal" style so no residual formatting diff is introduced.

$ws.Range("D2").Value = '22.394.80'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '1.568.63'
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.004'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3735'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.58%  '

$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3355'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07414'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.884'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.845'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.47%  '

$ws.Range("D16").Value = '1.562.68'
$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001109'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.86'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06671'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.142'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.11%  '

$ws.Range("D24").Value = '22.373.85'
$ws.Range("E24").Value = '  +0.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.365'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.523'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -9.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.998'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.83%  '

$ws.Range("D31").Value = '1.741.84'
$ws.Range("E31").Value = '  +0.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9989'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.974'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.905'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.636'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08377'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.367'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02444'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2239'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06374'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.353'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6170'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.794'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5762'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.046'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.41'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.218'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07294'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.36%  '
